# Update the "cryptos" price list (Coin / Link / Price / Volume(1h)) with
# freshly scraped values, per the GitHub Actions run on
# Fri Jun  2 10:34:01 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" values (column D) are plain digit strings (e.g. "1.001",
# "307.39") that Excel would otherwise auto-convert to numbers. Force those
# specific cells to text format first so the updated values stay strings,
# matching how the sheet was originally authored.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.102.14'
$ws.Range("D3").Value = '1.891.74'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '307.39'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D7").Value = '0.5138'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").Value = '0.3741'
$ws.Range("E8").Value = '  +3.23%  '
$ws.Range("D9").Value = '0.07213'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").Value = '0.9052'
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").Value = '1.898.23'
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = '94.95'
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("D15").Value = '5.272'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '0.000008479'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '14.46'
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '27.121.59'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '5.071'
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").Value = '2.128.20'
$ws.Range("E22").Value = '  +2.17%  '
$ws.Range("D23").Value = '10.56'
$ws.Range("E23").Value = '  +2.28%  '
$ws.Range("D24").Value = '6.405'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = '145.82'
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = '1.785'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '2.234'
$ws.Range("E27").Value = '  +8.69%  '
$ws.Range("D28").Value = '18.10'
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").Value = '114.61'
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").Value = '4.969'
$ws.Range("E30").Value = '  +6.21%  '
$ws.Range("D31").Value = '4.848'
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("D32").Value = '0.09188'
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").Value = '0.05093'
$ws.Range("D34").Value = '1.236'
$ws.Range("E34").Value = '  +7.54%  '
$ws.Range("D35").Value = '0.7703'
$ws.Range("E35").Value = '  +3.70%  '
$ws.Range("D36").Value = '2.992'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").Value = '3.291'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").Value = '2.631'
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").Value = '0.01997'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '0.5589'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '9.012'
$ws.Range("E42").Value = '  +5.74%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.657'
$ws.Range("E43").Value = '  +2.88%  '
$ws.Range("D44").Value = '117.73'
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").Value = '0.1511'
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("D46").Value = '0.4799'
$ws.Range("E46").Value = '  +1.83%  '
$ws.Range("D47").Value = '10.19'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = '0.9997'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '1.594'
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("D50").Value = '37.63'
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("E51").Value = '  +1.74%  '
